# Applies Peeyush's pass over Ashutosh's rows (2-13) in data_dict_reduced:
# fills in Data Type (I) / Type of analysis required (J) for each of those
# rows using the same "reviewed" orange-fill formatting already used further
# down the sheet (e.g. I26:K26), adjusts a couple of row heights so the
# wrapped text fits, updates the grade row's Data Type, drops the AutoFilter,
# and leaves the selection where Peeyush's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- copy the "reviewed" cell format (orange fill + border + wrap) onto
#     the I:J cells of every row being touched (K only for rows 2 & 3,
#     the rest keep their existing K formatting), using an already-styled
#     row (26) as the template so we reuse the existing style instead of
#     inventing a new one ---
$templateFormat = $ws.Range("I26:J26")
$templateFormat.Copy()
$targets = @(2,3,4,5,7,8,9,10,11,12,13)
foreach ($r in $targets) {
    $ws.Range("I" + $r + ":J" + $r).PasteSpecial(-4122) # xlPasteFormats
}
$excel.CutCopyMode = $false

# K2 and K3 also pick up the orange formatting (they stay empty)
$kFormat = $ws.Range("K26")
$kFormat.Copy()
$ws.Range("K2:K3").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- row 2: addr_state ---
$ws.Range("I2").Value = "Nominal"
$ws.Range("J2").Value = "rank-frequency plot`ngrouped bar chart by loan status"

# --- row 3: annual_inc ---
$ws.Range("I3").Value = "Interval"
$ws.Range("J3").Value = "histogram`nbox plot`nCorreleation with Loan Status`nbivariate analysis too"
$ws.Rows(3).RowHeight = 60

# --- row 4: collection_recovery_fee ---
$ws.Range("I4").Value = "Interval"
$ws.Range("J4").Value = "derived matrix of collection_recovery_fee and recoveries can help us understand the percentage the recovery fee is to the total recovery"
$ws.Rows(4).RowHeight = 60

# --- row 5: delinq_2yrs ---
$ws.Range("I5").Value = "Interval"
$ws.Range("J5").Value = "histogram`nbox plot`nCorreleation with Loan Status`nbivariate analysis too"
$ws.Rows(5).RowHeight = 60

# --- row 7: dti ---
$ws.Range("I7").Value = "Interval"
$ws.Range("J7").Value = "histogram`nbox plot`nCorreleation with Loan Status`nbivariate analysis too"

# --- row 8: earliest_cr_line ---
$ws.Range("I8").Value = "Date / Interval"
$ws.Range("J8").Value = "histogram`nbox plot`nCorreleation with Loan Status`nrelation between when the loan was asked (column name?)"
$ws.Rows(8).RowHeight = 75

# --- row 9: emp_length ---
$ws.Range("I9").Value = "Ordinal"
$ws.Range("J9").Value = "histogram`nbox plot`nCorreleation with Loan Status`nbivariate analysis too"
$ws.Rows(9).RowHeight = 60

# --- row 10: emp_title ---
$ws.Range("I10").Value = "Nominal"
$ws.Range("J10").Value = "Can be used to check if the loan status viz a viz a particular ogranisation is more or not. Is there are racket"
$ws.Rows(10).RowHeight = 45

# --- row 11: funded_amnt ---
$ws.Range("I11").Value = "Interval"
$ws.Range("J11").Value = "histogram`nbox plot`nCorreleation with Loan Status`nbivariate analysis too with funded_amnt_inv`nderived matrix funded_amnt_inv/funded_amnt"
$ws.Rows(11).RowHeight = 90

# --- row 12: funded_amnt_inv ---
$ws.Range("I12").Value = "Interval"
$ws.Range("J12").Value = "histogram`nbox plot`nCorreleation with Loan Status`nbivariate analysis too with funded_amnt`nderived matrix funded_amnt_inv/funded_amnt"
$ws.Rows(12).RowHeight = 90

# --- row 13: grade ---
$ws.Range("I13").Value = "Ordinal"
$ws.Range("J13").Value = "relation with Loan Status`nand other params`n"
$ws.Rows(13).RowHeight = 45

# --- remove the AutoFilter that used to sit over A1:K49 ---
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# --- leave the view where Peeyush's cursor ended up ---
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("J16").Select()
